$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7177
$ws.Range("C3").Value = 158128
$ws.Range("C4").Value = 149184
$ws.Range("C8").Value = 64
